$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "AppleI" (row 6) to "Apple-I" - same stats, just a label fix.
$ws.Range("A6").Value = "Apple-I"

# Insert a new row for "Apple-II" right after the (renamed) "Apple-I" row.
# This shifts the former rows 7-12 down to 8-13, carrying their formatting
# along with them (matches the row-insert Excel performed originally).
$ws.Rows.Item(7).Insert()

# Populate the new Apple-II row with its build stats.
$ws.Range("A7").Value = "Apple-II"
$ws.Range("B7").Value = 20211217
$ws.Range("C7").Value = "14,991 / 41,910"
$ws.Range("D7").Value = 0.36
$ws.Range("E7").Value = 17904
$ws.Range("F7").Value = "2,801,490 / 5,662,720"
$ws.Range("G7").Value = 0.49
$ws.Range("H7").Value = "362 / 553"
$ws.Range("I7").Value = 0.65
$ws.Range("J7").Value = "37 / 112"
$ws.Range("K7").Value = 0.33

# Reflect where the user's selection ended up after the edit.
$ws.Range("A14").Select()
